$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 11
$ws.Range("I2").Value = 12
$ws.Range("N2").Value = 3.55
$ws.Range("O2").Value = 1.37
$ws.Range("P2").Value = 1.85
$ws.Range("Q2").Value = 2.1
$ws.Range("U2").Value = 1.58
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AB2").Value = 6.4
$ws.Range("AC2").Value = 980
$ws.Range("AD2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("F3").Value = 2.9
$ws.Range("G3").Value = 2.92
$ws.Range("H3").Value = 2.72
$ws.Range("I3").Value = 2.74
$ws.Range("Z3").Value = 18
$ws.Range("AA3").Value = 42
$ws.Range("AB3").Value = 12
$ws.Range("AF3").Value = 19.5
$ws.Range("AH3").Value = 17.5
$ws.Range("AJ3").Value = 1000
$ws.Range("AL3").Value = 980
$ws.Range("AM3").Value = 1000
$ws.Range("P4").Value = 3.05
$ws.Range("Q4").Value = 1.46
$ws.Range("R4").Value = 1.83
$ws.Range("T4").Value = 1.69
$ws.Range("U4").Value = 2.38
$ws.Range("Z4").Value = 85
$ws.Range("AB4").Value = 13.5
$ws.Range("AC4").Value = 13.5
$ws.Range("AL4").Value = 25
$ws.Range("AN4").Value = 4.5
$ws.Range("J5").Value = 3.85
$ws.Range("K5").Value = 3.9
$ws.Range("P5").Value = 2.22
$ws.Range("T5").Value = 1.66
$ws.Range("U5").Value = 2.42
$ws.Range("X5").Value = 19.5
$ws.Range("Y5").Value = 1000
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 980
$ws.Range("AD5").Value = 1000
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 1000
$ws.Range("AN5").Value = 980
$ws.Range("AO5").Value = 1000
$ws.Range("F6").Value = 1.79
$ws.Range("G6").Value = 1.81
$ws.Range("I6").Value = 5.7
$ws.Range("K6").Value = 3.75
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 3.4
$ws.Range("O6").Value = 1.39
$ws.Range("P6").Value = 1.81
$ws.Range("Q6").Value = 2.18
$ws.Range("X6").Value = 13
$ws.Range("Y6").Value = 21
$ws.Range("Z6").Value = 1000
$ws.Range("AA6").Value = 160
$ws.Range("AC6").Value = 8.199999999999999
$ws.Range("AD6").Value = 1000
$ws.Range("AE6").Value = 1000
$ws.Range("AH6").Value = 24
$ws.Range("AI6").Value = 1000
$ws.Range("AJ6").Value = 1000
$ws.Range("AL6").Value = 1000
$ws.Range("AM6").Value = 130
$ws.Range("AN6").Value = 1000
$ws.Range("AO6").Value = 1000
$ws.Range("F7").Value = 3.25
$ws.Range("G7").Value = 3.3
$ws.Range("O7").Value = 1.33
$ws.Range("R7").Value = 1.36
$ws.Range("S7").Value = 3.65
$ws.Range("Y7").Value = 12
$ws.Range("Z7").Value = 16.5
$ws.Range("AA7").Value = 1000
$ws.Range("AB7").Value = 15.5
$ws.Range("AE7").Value = 1000
$ws.Range("AF7").Value = 24
$ws.Range("AH7").Value = 17.5
$ws.Range("AI7").Value = 1000
$ws.Range("AJ7").Value = 1000
$ws.Range("AL7").Value = 1000
$ws.Range("AN7").Value = 1000
$ws.Range("AO7").Value = 22
$ws.Range("G8").Value = 1.92
$ws.Range("N8").Value = 3.6
$ws.Range("O8").Value = 1.36
$ws.Range("R8").Value = 1.32
$ws.Range("T8").Value = 1.96
$ws.Range("U8").Value = 1.99
$ws.Range("X8").Value = 13.5
$ws.Range("Y8").Value = 15.5
$ws.Range("Z8").Value = 40
$ws.Range("AA8").Value = 970
$ws.Range("AD8").Value = 20
$ws.Range("AE8").Value = 1000
$ws.Range("AH8").Value = 23
$ws.Range("AI8").Value = 1000
$ws.Range("AL8").Value = 1000
$ws.Range("AM8").Value = 1000
$ws.Range("AN8").Value = 980
$ws.Range("AO8").Value = 1000
$ws.Range("F9").Value = 3.55
$ws.Range("H9").Value = 2.16
$ws.Range("I9").Value = 2.18
$ws.Range("Q9").Value = 1.86
$ws.Range("S9").Value = 3.2
$ws.Range("U9").Value = 2.3
$ws.Range("X9").Value = 18
$ws.Range("Z9").Value = 14.5
$ws.Range("AB9").Value = 1000
$ws.Range("AF9").Value = 1000
$ws.Range("AG9").Value = 18.5
$ws.Range("AJ9").Value = 1000
$ws.Range("AK9").Value = 1000
$ws.Range("AL9").Value = 1000
$ws.Range("AN9").Value = 1000
$ws.Range("AO9").Value = 16
$ws.Range("F10").Value = 2.54
$ws.Range("G10").Value = 2.56
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 3.4
$ws.Range("K10").Value = 3.25
$ws.Range("M10").Value = 1.11
$ws.Range("O10").Value = 1.44
$ws.Range("P10").Value = 1.69
$ws.Range("X10").Value = 10
$ws.Range("Y10").Value = 10.5
$ws.Range("Z10").Value = 1000
$ws.Range("AA10").Value = 1000
$ws.Range("AB10").Value = 10
$ws.Range("AD10").Value = 1000
$ws.Range("AE10").Value = 1000
$ws.Range("AF10").Value = 1000
$ws.Range("AG10").Value = 13
$ws.Range("AH10").Value = 1000
$ws.Range("AI10").Value = 1000
$ws.Range("AJ10").Value = 1000
$ws.Range("AK10").Value = 1000
$ws.Range("AL10").Value = 1000
$ws.Range("AM10").Value = 1000
$ws.Range("AN10").Value = 1000
$ws.Range("AO10").Value = 1000
$ws.Range("F11").Value = 2.38
$ws.Range("G11").Value = 2.4
$ws.Range("J11").Value = 3.2
$ws.Range("N11").Value = 2.98
$ws.Range("P11").Value = 1.65
$ws.Range("S11").Value = 4.8
$ws.Range("Y11").Value = 12.5
$ws.Range("Z11").Value = 1000
$ws.Range("AA11").Value = 1000
$ws.Range("AD11").Value = 1000
$ws.Range("AE11").Value = 1000
$ws.Range("AF11").Value = 16
$ws.Range("AH11").Value = 1000
$ws.Range("AI11").Value = 1000
$ws.Range("AJ11").Value = 1000
$ws.Range("AK11").Value = 1000
$ws.Range("AL11").Value = 1000
$ws.Range("AM11").Value = 1000
$ws.Range("AN11").Value = 1000
$ws.Range("AO11").Value = 1000
$ws.Range("H12").Value = 1.4
$ws.Range("I12").Value = 1.41
$ws.Range("P12").Value = 2.4
$ws.Range("Q12").Value = 1.67
$ws.Range("T12").Value = 1.95
$ws.Range("X12").Value = 25
$ws.Range("Y12").Value = 9.800000000000001
$ws.Range("Z12").Value = 9
$ws.Range("AA12").Value = 1000
$ws.Range("AB12").Value = 40
$ws.Range("AC12").Value = 1000
$ws.Range("AE12").Value = 1000
$ws.Range("AG12").Value = 1000
$ws.Range("AH12").Value = 30
$ws.Range("AI12").Value = 40
$ws.Range("AJ12").Value = 1000
$ws.Range("AK12").Value = 1000
$ws.Range("AM12").Value = 1000
$ws.Range("AN12").Value = 1000
$ws.Range("AO12").Value = 5.9
$ws.Range("F13").Value = 1.73
$ws.Range("G13").Value = 1.75
$ws.Range("N13").Value = 4.3
$ws.Range("O13").Value = 1.28
$ws.Range("Q13").Value = 1.86
$ws.Range("R13").Value = 1.43
$ws.Range("S13").Value = 3.15
$ws.Range("T13").Value = 1.85
$ws.Range("X13").Value = 17
$ws.Range("Z13").Value = 44
$ws.Range("AI13").Value = 90
$ws.Range("AN13").Value = 9.800000000000001
$ws.Range("AO13").Value = 1000
$ws.Range("H14").Value = 1.75
$ws.Range("I14").Value = 1.77
$ws.Range("Y14").Value = 9.800000000000001
$ws.Range("AO14").Value = 9.4
